$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet is being rebuilt: the old table had people as rows and phone
# brands as columns (with UserAvg / PhoneAvg / PhoneAvg-UserAvg helper rows
# and columns). The new table transposes that: brands as rows, people as
# columns, with a single BrandAvg column. Easiest/most reliable approach is
# to wipe all cell content and re-enter the new table from scratch.
# ---------------------------------------------------------------------------

$ws.Cells.ClearContents()

# ----- header row -----------------------------------------------------
$ws.Range("A1").Value = "Brand"
$ws.Range("B1").Value = "Olivia"
$ws.Range("C1").Value = "ZJY"
$ws.Range("D1").Value = "LHR"
$ws.Range("E1").Value = "ZHC"
$ws.Range("F1").Value = "ZYC"
$ws.Range("G1").Value = "Jason"
$ws.Range("H1").Value = "Tianen"
$ws.Range("I1").Value = "BrandAvg"

# ----- row 2 : Huawei ---------------------------------------------------
$ws.Range("A2").Value = "Huawei"
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 5

# ----- row 3 : Xiaomi ---------------------------------------------------
$ws.Range("A3").Value = "Xiaomi"
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 4
$ws.Range("H3").Value = 4

# ----- row 4 : Vivo -------------------------------------------------------
$ws.Range("A4").Value = "Vivo"
$ws.Range("E4").Value = 3.5

# ----- row 5 : Apple -------------------------------------------------------
$ws.Range("A5").Value = "Apple"
$ws.Range("B5").Value = 4.5
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 4.5
$ws.Range("F5").Value = 3
$ws.Range("H5").Value = 4.5

# ----- row 6 : Samsung ------------------------------------------------
$ws.Range("A6").Value = "Samsung"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 4
$ws.Range("E6").Value = 3.5
$ws.Range("F6").Value = 3
$ws.Range("H6").Value = 4

# ----- row 7 : Sony ---------------------------------------------------
$ws.Range("A7").Value = "Sony"
$ws.Range("E7").Value = 2.5

# ----- row 8 : Tianen ---------------------------------------------------
$ws.Range("A8").Value = "Tianen"
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 4
$ws.Range("E8").Value = 4.5
$ws.Range("F8").Value = 4
$ws.Range("H8").Value = 4.375

# ----- BrandAvg formulas (col I) ---------------------------------------
$ws.Range("I2").Formula = "= SUM(B2:H2)/COUNTIF(B2:H2, "">0"")"
$ws.Range("I3:I8").Formula = "= SUM(B3:H3)/COUNTIF(B3:H3, "">0"")"

# ----- column widths ------------------------------------------------------
$ws.Columns("A").ColumnWidth = 16.93
$ws.Columns("J").ColumnWidth = 17.6

# ----- selection ------------------------------------------------------
$ws.Range("J9").Select()
